$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.848.49"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "2.037.15"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.43"
$ws.Range("E5").Value = "  -1.41%  "
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  +3.40%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0818"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.66"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "2.338.60"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("D17").Value = "2.044.35"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "37.778.54"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.08"
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.79"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.14"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").Value = "  -2.27%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("E28").Value = "  -3.99%  "
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.29"
$ws.Range("E30").Value = "  -6.20%  "
$ws.Range("E31").Value = "  +1.49%  "
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("E33").Value = "  +4.30%  "
$ws.Range("E34").Value = "  -2.75%  "
$ws.Range("E35").Value = "  -2.52%  "
$ws.Range("E36").Value = "  +4.93%  "
$ws.Range("E37").Value = "  -5.40%  "
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "1.546.61"
$ws.Range("E40").Value = "  +4.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0218"
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.06"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.89"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0923"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("E47").Value = "  -4.52%  "
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.14"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").Value = "2.227.35"
$ws.Range("E51").Value = "  -1.14%  "
